# Regenerate merged AHB files
# Rename header labels from *_old/*_new to *_FV2304/*_FV2310, wrap the
# data range in an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- 1. Rename header cells (row 1) ---------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = ($oldHeaders[$i] -replace "_old$", "_FV2304")
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $cell = $ws.Cells.Item(1, 12 + $i)
    $cell.Value = ($newHeaders[$i] -replace "_new$", "_FV2310")
}

# --- 2. Wrap the used range in a Table -------------------------------------
$tableRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
